# Added ifo GDP component analysis preprocessing:
# diagonal shift of the qoq-error table (rows 45-53) by one column for
# the newly available evaluation horizon.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I45").Value = 0.4487415504340581
$ws.Range("J45").Value = 0.2388379152847414

$ws.Range("H46").Value = 0.5843816406042994
$ws.Range("I46").Value = 0.3744780054549828

$ws.Range("G47").Value = 0.3435754587486348
$ws.Range("H47").Value = 0.1336718235993181

$ws.Range("F48").Value = 0.2982442434965384
$ws.Range("G48").Value = 0.08834060834722172

$ws.Range("E49").Value = 0.2313828215604846
$ws.Range("F49").Value = 0.02147918641116785

$ws.Range("D50").Value = 0.201796619203768
$ws.Range("E50").Value = -0.00810701594554874

$ws.Range("C51").Value = 0.1836459624741271
$ws.Range("D51").Value = -0.02625767267518964

$ws.Range("B52").Value = 0.1656141382254278
$ws.Range("C52").Value = -0.04428949692388896

$ws.Range("B53").Value = -0.09587373626955231
